$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 14 de Junio de 2020 a las 20:14"

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 2151867
$ws.Cells.Item(4,3).Value = 9643
$ws.Cells.Item(4,4).Value = 857230
$ws.Cells.Item(4,5).Value = 1176961
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 149
$ws.Cells.Item(4,8).Value = 117676

# Row 7: India
$ws.Cells.Item(7,2).Value = 332739
$ws.Cells.Item(7,3).Value = 11113
$ws.Cells.Item(7,4).Value = 169465
$ws.Cells.Item(7,5).Value = 153760
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 315
$ws.Cells.Item(7,8).Value = 9514

# Row 9: España
$ws.Cells.Item(9,2).Value = 291008
$ws.Cells.Item(9,3).Value = 323
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(9,8).Value = 27136

# Row 12: Alemania
$ws.Cells.Item(12,2).Value = 187621
$ws.Cells.Item(12,3).Value = 198
$ws.Cells.Item(12,4).Value = 172200
$ws.Cells.Item(12,5).Value = 6553
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 1
$ws.Cells.Item(12,8).Value = 8868

# Row 14: Turquia
$ws.Cells.Item(14,2).Value = 178239
$ws.Cells.Item(14,3).Value = 1562
$ws.Cells.Item(14,4).Value = 151417
$ws.Cells.Item(14,5).Value = 22015
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 15
$ws.Cells.Item(14,8).Value = 4807

# Row 15: Chile
$ws.Cells.Item(15,2).Value = 174293
$ws.Cells.Item(15,3).Value = 6938
$ws.Cells.Item(15,4).Value = 137296
$ws.Cells.Item(15,5).Value = 33674
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 222
$ws.Cells.Item(15,8).Value = 3323

# Row 49: Israel
$ws.Cells.Item(49,2).Value = 19055
$ws.Cells.Item(49,3).Value = 83
$ws.Cells.Item(49,4).Value = 15375
$ws.Cells.Item(49,5).Value = 3380
$ws.Cells.Item(49,6).Value = 0
$ws.Cells.Item(49,7).Value = 0
$ws.Cells.Item(49,8).Value = 300

# Row 66: Marruecos
$ws.Cells.Item(66,2).Value = 8793
$ws.Cells.Item(66,3).Value = 101
$ws.Cells.Item(66,4).Value = 7765
$ws.Cells.Item(66,5).Value = 816
$ws.Cells.Item(66,6).Value = 0
$ws.Cells.Item(66,7).Value = 0
$ws.Cells.Item(66,8).Value = 212

# Row 101: Maldivas
$ws.Cells.Item(101,2).Value = 2035
$ws.Cells.Item(101,3).Value = 22
$ws.Cells.Item(101,4).Value = 1311
$ws.Cells.Item(101,5).Value = 716
$ws.Cells.Item(101,6).Value = 0
$ws.Cells.Item(101,7).Value = 0
$ws.Cells.Item(101,8).Value = 8

# Row 119: Paraguay
$ws.Cells.Item(119,1).Value = "Paraguay"
$ws.Cells.Item(119,2).Value = 1289
$ws.Cells.Item(119,3).Value = 28
$ws.Cells.Item(119,4).Value = 650
$ws.Cells.Item(119,5).Value = 628
$ws.Cells.Item(119,6).Value = 0
$ws.Cells.Item(119,7).Value = 0
$ws.Cells.Item(119,8).Value = 11

# Row 120: Madagascar
$ws.Cells.Item(120,1).Value = "Madagascar"
$ws.Cells.Item(120,2).Value = 1272
$ws.Cells.Item(120,3).Value = 20
$ws.Cells.Item(120,4).Value = 367
$ws.Cells.Item(120,5).Value = 895
$ws.Cells.Item(120,6).Value = 0
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = 10

# Row 144: Malaui
$ws.Cells.Item(144,1).Value = "Malaui"
$ws.Cells.Item(144,2).Value = 547
$ws.Cells.Item(144,3).Value = 18
$ws.Cells.Item(144,4).Value = 69
$ws.Cells.Item(144,5).Value = 472
$ws.Cells.Item(144,6).Value = 0
$ws.Cells.Item(144,7).Value = 1
$ws.Cells.Item(144,8).Value = 6

# Row 145: Ruanda
$ws.Cells.Item(145,1).Value = "Ruanda"
$ws.Cells.Item(145,2).Value = 541
$ws.Cells.Item(145,3).Value = 0
$ws.Cells.Item(145,4).Value = 332
$ws.Cells.Item(145,5).Value = 207
$ws.Cells.Item(145,6).Value = 0
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 2

# Row 146: Togo
$ws.Cells.Item(146,1).Value = "Togo"
$ws.Cells.Item(146,2).Value = 530
$ws.Cells.Item(146,3).Value = 0
$ws.Cells.Item(146,4).Value = 291
$ws.Cells.Item(146,5).Value = 226
$ws.Cells.Item(146,6).Value = 0
$ws.Cells.Item(146,7).Value = 0
$ws.Cells.Item(146,8).Value = 13

# Row 149: Estado de Palestina
$ws.Cells.Item(149,2).Value = 492
$ws.Cells.Item(149,3).Value = 3
$ws.Cells.Item(149,4).Value = 415
$ws.Cells.Item(149,5).Value = 74
$ws.Cells.Item(149,6).Value = 0
$ws.Cells.Item(149,7).Value = 0
$ws.Cells.Item(149,8).Value = 3

# Row 150: Suazilandia
$ws.Cells.Item(150,2).Value = 490
$ws.Cells.Item(150,3).Value = 4
$ws.Cells.Item(150,4).Value = 249
$ws.Cells.Item(150,5).Value = 237
$ws.Cells.Item(150,6).Value = 0
$ws.Cells.Item(150,7).Value = 1
$ws.Cells.Item(150,8).Value = 4

# Row 206: Islas Malvinas
$ws.Cells.Item(206,1).Value = "Islas Malvinas"
$ws.Cells.Item(206,2).Value = 13
$ws.Cells.Item(206,3).Value = 0
$ws.Cells.Item(206,4).Value = 13
$ws.Cells.Item(206,5).Value = 0
$ws.Cells.Item(206,6).Value = 0
$ws.Cells.Item(206,7).Value = 0
$ws.Cells.Item(206,8).Value = 0

# Row 207: Groenlandia
$ws.Cells.Item(207,1).Value = "Groenlandia"
$ws.Cells.Item(207,2).Value = 13
$ws.Cells.Item(207,3).Value = 0
$ws.Cells.Item(207,4).Value = 13
$ws.Cells.Item(207,5).Value = 0
$ws.Cells.Item(207,6).Value = 0
$ws.Cells.Item(207,7).Value = 0
$ws.Cells.Item(207,8).Value = 0

# Row 208: Santa Sede
$ws.Cells.Item(208,1).Value = "Santa Sede"
$ws.Cells.Item(208,2).Value = 12
$ws.Cells.Item(208,3).Value = 0
$ws.Cells.Item(208,4).Value = 12
$ws.Cells.Item(208,5).Value = 0
$ws.Cells.Item(208,6).Value = 0
$ws.Cells.Item(208,7).Value = 0
$ws.Cells.Item(208,8).Value = 0

# Row 209: Islas Turcas y Caicos
$ws.Cells.Item(209,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209,2).Value = 12
$ws.Cells.Item(209,3).Value = 0
$ws.Cells.Item(209,4).Value = 11
$ws.Cells.Item(209,5).Value = 0
$ws.Cells.Item(209,6).Value = 0
$ws.Cells.Item(209,7).Value = 0
$ws.Cells.Item(209,8).Value = 1
